# "Add PPL data and targets" - one more row of portfolio company data
# (PPL Corp.) added to the ITR input/target sheets and the Portfolio sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ITR input data" - new company row 33 (PPL Corp.)
# ---------------------------------------------------------------------
$wsItr = $wb.Worksheets.Item("ITR input data")

# Row 32 carries the formatting we want to replicate onto row 33 for the
# cells whose style needs to change (O33/P33 drop the "blank row" font,
# AE33:AI33 pick up the highlighted-column fill used elsewhere in the
# table). Copy formats only so we don't disturb values we set below.
$wsItr.Range("O32:P32").Copy() | Out-Null
$wsItr.Range("O33:P33").PasteSpecial(-4122) | Out-Null
$wsItr.Range("AE32:AI32").Copy() | Out-Null
$wsItr.Range("AE33:AI33").PasteSpecial(-4122) | Out-Null

$wsItr.Range("A33").Value = "PPL Corp."
$wsItr.Range("B33").Value = "9N3UAJSNOUXFKQLF3V18"
$wsItr.Range("C33").Value = "US69351T1060"
$wsItr.Range("D33").Value = "US"
$wsItr.Range("E33").Value = "North America"
$wsItr.Range("F33").Value = "Electricity Utilities"
$wsItr.Range("G33").Value = "equity"
$wsItr.Range("H33").Value = "USD"
$wsItr.Range("I33").Value = 44196
$wsItr.Range("J33").Value = 19865342074
$wsItr.Range("K33").Value = 7769000000
$wsItr.Range("L33").Value = 40943342074
$wsItr.Range("M33").Value = 41758342074
$wsItr.Range("N33").Value = 45680000000
$wsItr.Range("O33").Value = "Mt CO2"
$wsItr.Range("P33").Value = "TWh"

$wsItr.Range("AE33").Value = 30.088487220000001
$wsItr.Range("AF33").Value = 30.24837145
$wsItr.Range("AG33").Value = 31.611469039999999
$wsItr.Range("AH33").Value = 28.778915319999999
$wsItr.Range("AI33").Value = 28.07780713

$wsItr.Range("AS33").Value = 38.355258640000002
$wsItr.Range("AT33").Value = 37.442832350000003
$wsItr.Range("AU33").Value = 39.590075179999999
$wsItr.Range("AV33").Value = 35.152931719999998
$wsItr.Range("AW33").Value = 32.487984334642732

$wsItr.Range("A33:C33").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "ITR target input data" - two new target rows (42, 43) for PPL
# ---------------------------------------------------------------------
$wsTgt = $wb.Worksheets.Item("ITR target input data")

# Row 41 has the formatting we want for the new rows; I42/I43 already
# carry their own (number-format) style so re-apply it after the
# row-level paste overwrites it.
$wsTgt.Range("A41:L41").Copy() | Out-Null
$wsTgt.Range("A42:L42").PasteSpecial(-4122) | Out-Null
$wsTgt.Range("A43:L43").PasteSpecial(-4122) | Out-Null
$wsTgt.Range("I44").Copy() | Out-Null
$wsTgt.Range("I42:I43").PasteSpecial(-4122) | Out-Null

$wsTgt.Range("A42").Value = "PPL Corp."
$wsTgt.Range("B42").Value = "9N3UAJSNOUXFKQLF3V18"
$wsTgt.Range("C42").Value = "US69351T1060"
$wsTgt.Range("D42").Value = 2050
$wsTgt.Range("E42").Value = "absolute"
$wsTgt.Range("F42").Value = "S1+S2"
$wsTgt.Range("G42").Value = 2021
$wsTgt.Range("H42").Value = 2010
$wsTgt.Range("I42").Formula = "=60736086+1597157"
$wsTgt.Range("J42").Value = "t CO2"
$wsTgt.Range("K42").Value = 2035
$wsTgt.Range("L42").Value = 0.7

$wsTgt.Range("A43").Value = "PPL Corp."
$wsTgt.Range("B43").Value = "9N3UAJSNOUXFKQLF3V18"
$wsTgt.Range("C43").Value = "US69351T1060"
$wsTgt.Range("D43").Value = 2050
$wsTgt.Range("E43").Value = "absolute"
$wsTgt.Range("F43").Value = "S1+S2"
$wsTgt.Range("G43").Value = 2021
$wsTgt.Range("H43").Value = 2010
$wsTgt.Range("I43").Formula = "=60736086+1597157"
$wsTgt.Range("J43").Value = "t CO2"
$wsTgt.Range("K43").Value = 2040
$wsTgt.Range("L43").Value = 0.8

$wsTgt.Range("A43").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "ITR input data (2)" - view only changes (selection moved to
# the full row 40)
# ---------------------------------------------------------------------
$wsItr2 = $wb.Worksheets.Item("ITR input data (2)")
$wsItr2.Range("A40:XFD40").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Portfolio" - new row 33 (PPL Corp. again, random benchmark
# value). This is the sheet left active/selected at the end.
# ---------------------------------------------------------------------
$wsPort = $wb.Worksheets.Item("Portfolio")

$wsPort.Range("A32:D32").Copy() | Out-Null
$wsPort.Range("A33:D33").PasteSpecial(-4122) | Out-Null

$wsPort.Range("A33").Value = "PPL Corp."
$wsPort.Range("B33").Value = "9N3UAJSNOUXFKQLF3V18"
$wsPort.Range("C33").Value = "US69351T1060"
$wsPort.Range("D33").Value = "US69351T1060"
$wsPort.Range("E33").Formula = "=RANDBETWEEN(35000,250000)"

$wsPort.Range("E33").Select() | Out-Null

Write-Host "Added PPL Corp. row + targets"
